$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New SVR parameter headers (row 1)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New SVR parameter values (row 2)
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0.05
$ws.Range("M2").Value = 20

$ws.Range("I8").Select()
